$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 74 values (revised figures)
$ws.Range("B74").Value = 33.3
$ws.Range("C74").Value = 77.8
$ws.Range("D74").Value = 130.1
$ws.Range("E74").Value = 121.8

# Add new row 75 (01-04-2021)
$ws.Range("A75").Value = "'01-04-2021"
$ws.Range("A75").Style = "Normal"
$ws.Range("B75").Value = 36.5
$ws.Range("C75").Value = 76.7
$ws.Range("D75").Value = 129
$ws.Range("E75").Value = 117.6
